$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-ambiguous string assignments
$ws.Range("D2").Value = "59.431.57"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "2.526.42"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("E6").Value = "  -2.87%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("D9").Value = "2.530.28"
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("E10").Value = "  -0.12%  "
$ws.Range("E11").Value = "  +1.30%  "
$ws.Range("E12").Value = "  -2.27%  "
$ws.Range("E13").Value = "  +1.15%  "
$ws.Range("D14").Value = "2.972.47"
$ws.Range("E15").Value = "  -2.31%  "
$ws.Range("D16").Value = "59.363.97"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("D18").Value = "2.506.50"
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("E19").Value = "  -3.00%  "
$ws.Range("E20").Value = "  -1.34%  "
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("E25").Value = "  -3.74%  "
$ws.Range("E26").Value = "  +1.65%  "
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("E29").Value = "  -0.67%  "
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("E32").Value = "  +1.40%  "
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("E34").Value = "  +1.29%  "
$ws.Range("E35").Value = "  -7.03%  "
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("E37").Value = "  -3.95%  "
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("E42").Value = "  -5.57%  "
$ws.Range("E43").Value = "  -7.29%  "
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("E47").Value = "  -1.90%  "
$ws.Range("E48").Value = "  -0.77%  "
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("E51").Value = "  -2.03%  "

# Numeric-looking strings that must remain text: force text format, then reset style
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "536.47"
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "139.77"
$cell.Style = "Normal"
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "322.06"
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "61.95"
$cell.Style = "Normal"
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.423"
$cell.Style = "Normal"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "160.19"
$cell.Style = "Normal"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "4.21"
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "286.09"
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "5.25"
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.997"
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.598"
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "10.85"
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "124.06"
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "18.57"
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.0511"
$cell.Style = "Normal"
